$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '39.524.93'
$ws.Range("E2").Value = '  -3.39%  '
$ws.Range("D3").Value = '2.297.34'
$ws.Range("E3").Value = '  -4.56%  '
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '304.82'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.30%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '81.01'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -8.07%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.514'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -4.08%  '
$ws.Range("E8").Value = '  +0.20%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.467'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.16%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0780'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.15%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '28.34'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -10.10%  '
$ws.Range("E12").Value = '  -0.43%  '
$ws.Range("D13").Value = '2.653.31'
$ws.Range("E13").Value = '  -4.34%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.15'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -7.57%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.34'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -7.72%  '
$ws.Range("D16").Value = '2.310.06'
$ws.Range("E16").Value = '  -4.45%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.728'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.01%  '
$ws.Range("D18").Value = '39.466.55'
$ws.Range("E18").Value = '  -3.14%  '
$ws.Range("D19").Value = '0.0₃0873'
$ws.Range("E19").Value = '  -4.71%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.87'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.56%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '66.74'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -6.40%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.06'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -7.25%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '229.73'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.05%  '
$ws.Range("E24").Value = '  -0.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.44'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -8.47%  '
$ws.Range("E26").Value = '  -4.30%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.57'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.64%  '
$ws.Range("E28").Value = '  -1.85%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.95'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.59%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '149.77'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.59%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '31.96'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.06%  '
$ws.Range("E32").Value = '  +0.13%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.89'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.30%  '
$ws.Range("E34").Value = '  -0.74%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0695'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.34%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.111'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.89%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.66'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -7.57%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0959'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.83%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '14.99'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -8.19%  '
$ws.Range("E40").Value = '  -7.58%  '
$ws.Range("E41").Value = '  -4.89%  '
$ws.Range("E42").Value = '  -1.59%  '
$ws.Range("D43").Value = '1.946.03'
$ws.Range("E43").Value = '  -1.72%  '
$ws.Range("E44").Value = '  -6.13%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.60'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -9.84%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.15'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.95%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.57'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -9.60%  '
$ws.Range("D48").Value = '2.523.91'
$ws.Range("E48").Value = '  -4.46%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '89.75'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.20%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '67.56'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -7.29%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '48.00'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.90%  '
